$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.936.23"
$ws.Range("E2").Value = "  -3.37%  "

$ws.Range("D3").Value = "1.858.73"
$ws.Range("E3").Value = "  -2.59%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.79"
$ws.Range("E5").Value = "  -2.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4360"
$ws.Range("E7").Value = "  -4.87%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3699"
$ws.Range("E8").Value = "  -2.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07486"
$ws.Range("E9").Value = "  -2.95%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9399"
$ws.Range("E10").Value = "  -4.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.27"
$ws.Range("E11").Value = "  -3.56%  "

$ws.Range("D12").Value = "1.840.04"
$ws.Range("E12").Value = "  -3.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.718"
$ws.Range("E13").Value = "  -3.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.435"
$ws.Range("E14").Value = "  -4.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06868"
$ws.Range("E15").Value = "  -2.51%  "

$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.59"
$ws.Range("E17").Value = "  -2.55%  "

$ws.Range("E18").Value = "  -4.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.94"
$ws.Range("E20").Value = "  -4.03%  "

$ws.Range("D21").Value = "27.907.87"
$ws.Range("E21").Value = "  -3.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.112"
$ws.Range("E22").Value = "  -3.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.02"
$ws.Range("E23").Value = "  +1.25%  "

$ws.Range("D24").Value = "2.104.87"
$ws.Range("E24").Value = "  -1.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.009"
$ws.Range("E25").Value = "  -4.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.52"
$ws.Range("E26").Value = "  -2.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.36"
$ws.Range("E27").Value = "  -3.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.383"
$ws.Range("E28").Value = "  -4.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.49"
$ws.Range("E29").Value = "  -3.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.732"
$ws.Range("E30").Value = "  -7.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08983"
$ws.Range("E31").Value = "  -3.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8044"
$ws.Range("E32").Value = "  -6.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.831"
$ws.Range("E33").Value = "  -4.79%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.172"
$ws.Range("E34").Value = "  -6.15%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.956"
$ws.Range("E35").Value = "  -4.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.001"
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05474"
$ws.Range("E37").Value = "  -4.14%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.111"
$ws.Range("E38").Value = "  -3.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01980"
$ws.Range("E39").Value = "  -2.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.941"
$ws.Range("E40").Value = "  +1.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5253"
$ws.Range("E41").Value = "  -4.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.004"
$ws.Range("E42").Value = "  -5.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1681"
$ws.Range("E43").Value = "  -4.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.783"
$ws.Range("E44").Value = "  -5.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06774"
$ws.Range("E45").Value = "  -1.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4902"
$ws.Range("E46").Value = "  -5.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.67"
$ws.Range("E47").Value = "  -4.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.99"
$ws.Range("E48").Value = "  -3.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.919"
$ws.Range("E49").Value = "  -9.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.680"
$ws.Range("E50").Value = "  -5.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.000"
$ws.Range("E51").Value = "  -0.17%  "
